$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 320 (weekly update: new price observation),
# pushing the existing rows 320-366 down to 321-367.
$ws.Rows(320).Insert()

$ws.Cells.Item(320, 1).Value  = 1
$ws.Cells.Item(320, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(320, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(320, 4).Value  = 44984
$ws.Cells.Item(320, 5).Value  = 15
$ws.Cells.Item(320, 6).Value  = "Fruta"
$ws.Cells.Item(320, 7).Value  = 100108
$ws.Cells.Item(320, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(320, 9).Value  = 100108006
$ws.Cells.Item(320, 10).Value = "Plátano"
$ws.Cells.Item(320, 11).Value = "Sin especificar"
$ws.Cells.Item(320, 12).Value = "Pintón"
$ws.Cells.Item(320, 13).Value = 120
$ws.Cells.Item(320, 14).Value = 19000
$ws.Cells.Item(320, 15).Value = 20000
$ws.Cells.Item(320, 16).Value = 19500
$ws.Cells.Item(320, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(320, 18).Value = "Ecuador"
$ws.Cells.Item(320, 19).Value = 975
$ws.Cells.Item(320, 20).Value = 20
